$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Register the two new cell styles used by the taller MSG_012 row further
# down (left+wrap for the memo column, center+vertical-center for the
# STT/Code columns). They are built transiently on two already-formatted
# cells and then reverted so those cells keep their original look -
# Excel keeps the now-unused style records in the style table, which is
# exactly how the target workbook ends up with extra, pre-registered
# cellXfs entries (plus a small leftover 8pt font picked up along the way).
$ws.Range("E13").Phonetics.Font.Size = 8
$ws.Range("E13").WrapText = $true
$ws.Range("E13").Font.Size = 11
$ws.Range("E13").WrapText = $false

$ws.Range("C14").Phonetics.Font.Size = 8
$ws.Range("C14").VerticalAlignment = -4108
$ws.Range("C14").Font.Size = 11
$ws.Range("C14").VerticalAlignment = -4107

# ---------------------------------------------------------------------------
# Login / logout / change-password message list updates.
#
# Row 14 is filled in with the new MSG_011 "logout confirmation" message.
$ws.Range("C14").Value = 11
$ws.Range("D14").Value = "MSG_011"
$ws.Range("E14").Value = "Are you sure you want to logout?"

# Row 15 is filled in with the new MSG_012 "session expired" message.
$ws.Range("C15").Value = 12
$ws.Range("D15").Value = "MSG_012"
$ws.Range("E15").Value = "Sesion Expired" + [char]10 + "Please login again."

# The new message is two lines long, so the row is made taller.
$ws.Rows(15).RowHeight = 30

# E15 gets the left aligned, word-wrapped style ...
$ws.Range("E15").HorizontalAlignment = -4131
$ws.Range("E15").WrapText = $true

# ... while C15:D15 get the centered (horizontal + vertical) style.
$ws.Range("C15:D15").HorizontalAlignment = -4108
$ws.Range("C15:D15").VerticalAlignment = -4108

# Row 13 (MSG_010): the memo text is updated last, from the old password
# complexity rule text to the new, simpler "Invalid password." message.
$ws.Range("E13").Value = "Invalid password."

# ---------------------------------------------------------------------------
# Sheet-level view tweaks.
$ws.Range("E7").Select()
